$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# 1. "Certified " + "SAFe" + " 5 Agilist" -> single run "Certified SAFe 5 Agilist"
Replace-Text "Certified SAFe 5 Agilist" "Certified SAFe 5 Agilist"

# 2. "Working Knowledge ... HP " + "QC,HP" + " UFT, ... HP PC," -> single run
Replace-Text "Working Knowledge of various tools like HP QC,HP UFT, HP ALM,HP AGM, Microfocus ALM Octane, HP PC," "Working Knowledge of various tools like HP QC,HP UFT, HP ALM,HP AGM, Microfocus ALM Octane, HP PC,"

# 3. "Kingdom(" -> "Kingdom (" (keep as its own run, between "Atos Syntel United " and "Aug 2020)")
Replace-Text "Kingdom(" "Kingdom ("

# 4. "- " + "ePDSM" + ", Microfocus Octane, ... Prometheus," -> single run
Replace-Text "- ePDSM, Microfocus Octane, Microfocus ALM, GITLAB, Grafana, Prometheus," "- ePDSM, Microfocus Octane, Microfocus ALM, GITLAB, Grafana, Prometheus,"

# 5. "Sharepoint, Jenkins, " + "BlueOcean" -> single run
Replace-Text "Sharepoint, Jenkins, BlueOcean" "Sharepoint, Jenkins, BlueOcean"

# 6. "BDD ," + " TDD methods" -> single run
Replace-Text "BDD , TDD methods" "BDD , TDD methods"

# 7. "How do you assign story " + "points" -> single run
Replace-Text "How do you assign story points" "How do you assign story points"

# 8. "Tracking charts as Burn down / burn " + "up ," + " velocity tracking . What are they and how do you use " + "them" -> single run
Replace-Text "Tracking charts as Burn down / burn up , velocity tracking . What are they and how do you use them" "Tracking charts as Burn down / burn up , velocity tracking . What are they and how do you use them"

# 9. "How would you manage virtual meetings due to work from " + "home." + " Answer " + "–" + " Try to have video calls often" -> single run
Replace-Text "How would you manage virtual meetings due to work from home. Answer – Try to have video calls often" "How would you manage virtual meetings due to work from home. Answer – Try to have video calls often"

# 10. "Team size and roles. What is the ideal size of an Agile " + "team" -> single run
Replace-Text "Team size and roles. What is the ideal size of an Agile team" "Team size and roles. What is the ideal size of an Agile team"

# 11. "Where do you see yourself in 5 " + "years" -> single run
Replace-Text "Where do you see yourself in 5 years" "Where do you see yourself in 5 years"
